$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.505.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.493.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.20%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.66'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.647'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.99'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000307'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.92%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.50'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.046.80'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '602.37'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.513.71'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.60'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.486.13'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.96%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.985'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.21'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.58'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.12'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.64'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.05'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.07%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.43'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.18'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +15.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.43'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.47%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.25%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.71'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '518.05'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.396'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.599.31'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.76'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.85%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0462'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.93'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.77'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.36%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.00'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.36'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -9.17%  '
